$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 'FUNC_BaiHoc_TC1'
$ws.Range("B8").Value = 'Chọn bài học và hiển thị nội dung bài học'
$ws.Range("C8").Value = 'Không có'
$ws.Range("D8").Value = 'Không có'
$ws.Range("E8").Value = 'Chọn tuần 1 . Chọn bài học tập đọc . Chọn bài Mẹ Ốm. '
$ws.Range("F8").Value = 'Khi chọn tuần 1 , màn hình hiện ra các bài học thuộc tuần 1 : tập đọc , tập làm văn , luyện từ và câu , chính tả , kể chuyện. Sau khi chọn tập đọc thì màn hình sẽ hiện ra hai bài tập đọc : Dế mèn bênh vực kẻ yếu và bài Mẹ Ốm. Sau khi click chọn bài Mẹ Ốm thì màn hình sẽ hiện ra nội dung bài Mẹ Ốm: " Mọi hôm mẹ thích vui chơi..."'
$ws.Rows.Item(8).RowHeight = 102

$ws.Range("A9").Value = 'FUNC_BaiHoc_TC2'
$ws.Range("B9").Value = 'Trở lại màn hình trước để chọn bài học khác'
$ws.Range("C9").Value = 'Không có'
$ws.Range("D9").Value = 'Không có'
$ws.Range("E9").Value = 'Chọn tuần 1 . Chọn bài học tập đọc . Trở lai màn hình trước đó. Chọn bài học chính tả.'
$ws.Range("F9").Value = 'Khi chọn tuần 1 , màn hình hiện ra các bài học thuộc tuần 1 : tập đọc , tập làm văn , luyện từ và câu , chính tả , kể chuyện. Sau khi chọn tập đọc thì màn hình sẽ hiện ra hai bài tập đọc : Dế mèn bênh vực kẻ yếu và bài Mẹ Ốm. Trở về màn hình trước đó , chọn bài học chính tả.'
$ws.Rows.Item(9).RowHeight = 89.25

$ws.Range("A10").Value = 'FUNC_BaiHoc_TC3'
$ws.Range("B10").Value = 'Trở lại màn hình chính của phần mềm'
$ws.Range("C10").Value = 'Không có'
$ws.Range("D10").Value = 'Không có'
$ws.Range("E10").Value = 'Chọn tuần học 1 . Chọn bài học chính tả. Chọn bài mẹ ốm . Trở lại màn hính chính.'
$ws.Range("F10").Value = 'Về màn hình chính của phần mềm: gồm 2 phần học tập và giải trí bất cứ lúc nào'
$ws.Rows.Item(10).RowHeight = 38.25


# Touch rows 11-12 (A and H only) so they materialize as styled-but-empty cells
$ws.Range("A11").Font.Name = $ws.Range("A11").Font.Name
$ws.Range("H11").Font.Name = $ws.Range("H11").Font.Name
$ws.Range("A12").Font.Name = $ws.Range("A12").Font.Name
$ws.Range("H12").Font.Name = $ws.Range("H12").Font.Name

# Touch row 13 (all columns) - bottom boundary row of the table
$ws.Range("A13:H13").Font.Name = "Times New Roman"

# Convert the data range into an Excel Table (ListObject)
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:H13"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight1"

# Update view: scroll to show the new rows, select A10
$ws.Range("A10").Select()
$ws.Application.ActiveWindow.ScrollRow = 8
